# erdri_cds.xlsx: rename header/column labels and tidy up the "code" suffixes
# used in the data_types (column D) hints, per the source notebook update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): relabel the columns -----------------------------
$ws.Range("A1").Value = "data_model_section"
$ws.Range("B1").Value = "data_field_name"
$ws.Range("C1").Value = "description"
$ws.Range("D1").Value = "data_types"
$ws.Range("E1").Value = "required"
$ws.Range("F1").Value = "comment"

# --- Column D "data_types" hints: drop the redundant "code"/"codes" word -
$ws.Range("D10").Value = "orpha, alpha, icd-9, icd-9-cm, icd-10"
$ws.Range("D11").Value = "hgvs, hgnc, omim"
$ws.Range("D12").Value = "hpo, hgvs"

# --- Column widths: A/B widened (best-fit) for the longer header/labels --
$ws.Columns.Item(1).ColumnWidth = 19.3
$ws.Columns.Item(2).ColumnWidth = 32.15

# --- Selection moves to F1 -------------------------------------------------
$ws.Range("F1").Select()
